$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.464.81'
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('D3').Value = '3.380.32'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '''574.98'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.45%  '
$ws.Range('D6').Value = '''137.39'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.65%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '3.379.12'
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').Value = '''0.474'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.16%  '
$ws.Range('E10').Value = '  -1.69%  '
$ws.Range('E11').Value = '  +1.01%  '
$ws.Range('E12').Value = '  -0.94%  '
$ws.Range('D13').Value = '3.956.74'
$ws.Range('E13').Value = '  +0.01%  '
$ws.Range('E14').Value = '  +2.52%  '
$ws.Range('D15').Value = '''0.0000175'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.11%  '
$ws.Range('D16').Value = '''26.20'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.77%  '
$ws.Range('D17').Value = '3.374.36'
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('D18').Value = '61.573.39'
$ws.Range('E18').Value = '  +0.21%  '
$ws.Range('D19').Value = '''14.05'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.52%  '
$ws.Range('E20').Value = '  +0.76%  '
$ws.Range('E21').Value = '  -1.21%  '
$ws.Range('D22').Value = '''377.28'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.32%  '
$ws.Range('E23').Value = '  -3.60%  '
$ws.Range('D24').Value = '3.511.12'
$ws.Range('E24').Value = '  -0.13%  '
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('E26').Value = '  +5.87%  '
$ws.Range('D27').Value = '''71.64'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.06%  '
$ws.Range('D28').Value = '''1.72'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.75%  '
$ws.Range('D29').Value = '''7.53'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.14%  '
$ws.Range('D30').Value = '''0.995'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.37%  '
$ws.Range('E31').Value = '  +1.47%  '
$ws.Range('D32').Value = '''0.161'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.08%  '
$ws.Range('E33').Value = '  +1.08%  '
$ws.Range('E34').Value = '  +0.08%  '
$ws.Range('D35').Value = '''23.68'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.73%  '
$ws.Range('E36').Value = '  -5.29%  '
$ws.Range('D37').Value = '''6.83'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.27%  '
$ws.Range('E38').Value = '  -1.26%  '
$ws.Range('D39').Value = '''166.01'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.78%  '
$ws.Range('D40').Value = '''0.0774'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.91%  '
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('E42').Value = '  +3.22%  '
$ws.Range('D43').Value = '''0.775'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.69%  '
$ws.Range('D44').Value = '''41.64'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.21%  '
$ws.Range('E45').Value = '  -0.60%  '
$ws.Range('E46').Value = '  -0.07%  '
$ws.Range('D47').Value = '''24.61'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.25%  '
$ws.Range('E48').Value = '  -1.77%  '
$ws.Range('D49').Value = '''22.65'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.46%  '
$ws.Range('D50').Value = '2.366.62'
$ws.Range('E50').Value = '  +3.81%  '
$ws.Range('E51').Value = '  -1.42%  '
